$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.44654466666667
$ws.Range("H2").Value = 58.339634
$ws.Range("I2").Value = 0.7934109702307454
$ws.Range("J2").Value = 0.7934109702307454
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 102.8289443333334
$ws.Range("N2").Value = 308.486833
$ws.Range("O2").Value = 0.5559120396302444
$ws.Range("P2").Value = 0.5559120396302443
$ws.Range("Q2").Value = 1999.667659004347
$ws.Range("R2").Value = 17997.00893103913
$ws.Range("S2").Value = 0.4410667107259847
$ws.Range("T2").Value = 0.4410667107259847

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.44654466666667
$ws.Range("H3").Value = 58.339634
$ws.Range("I3").Value = 0.7934109702307454
$ws.Range("J3").Value = 0.7934109702307454
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 63.66262833333334
$ws.Range("N3").Value = 190.987885
$ws.Range("O3").Value = 0.3441717873742006
$ws.Range("P3").Value = 0.3441717873742006
$ws.Range("Q3").Value = 1238.018145481566
$ws.Range("R3").Value = 11142.16330933409
$ws.Range("S3").Value = 0.2730696717466143
$ws.Range("T3").Value = 0.2730696717466143

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.44654466666667
$ws.Range("H4").Value = 58.339634
$ws.Range("I4").Value = 0.7934109702307454
$ws.Range("J4").Value = 0.7934109702307454
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.481835
$ws.Range("N4").Value = 55.445505
$ws.Range("O4").Value = 0.09991617299555507
$ws.Range("P4").Value = 0.09991617299555505
$ws.Range("Q4").Value = 359.4078298494634
$ws.Range("R4").Value = 3234.67046864517
$ws.Range("S4").Value = 0.07927458775814634
$ws.Range("T4").Value = 0.07927458775814633

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.524415666666667
$ws.Range("H5").Value = 7.573247
$ws.Range("I5").Value = 0.1029951139231878
$ws.Range("J5").Value = 0.1029951139231878
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 102.8289443333334
$ws.Range("N5").Value = 308.486833
$ws.Range("O5").Value = 0.5559120396302444
$ws.Range("P5").Value = 0.5559120396302443
$ws.Range("Q5").Value = 259.5829980618612
$ws.Range("R5").Value = 2336.246982556751
$ws.Range("S5").Value = 0.05725622385298872
$ws.Range("T5").Value = 0.05725622385298871

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.524415666666667
$ws.Range("H6").Value = 7.573247
$ws.Range("I6").Value = 0.1029951139231878
$ws.Range("J6").Value = 0.1029951139231878
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 63.66262833333334
$ws.Range("N6").Value = 190.987885
$ws.Range("O6").Value = 0.3441717873742006
$ws.Range("P6").Value = 0.3441717873742006
$ws.Range("Q6").Value = 160.7109363458439
$ws.Range("R6").Value = 1446.398427112595
$ws.Range("S6").Value = 0.03544801244975297
$ws.Range("T6").Value = 0.03544801244975297

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.524415666666667
$ws.Range("H7").Value = 7.573247
$ws.Range("I7").Value = 0.1029951139231878
$ws.Range("J7").Value = 0.1029951139231878
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.481835
$ws.Range("N7").Value = 55.445505
$ws.Range("O7").Value = 0.09991617299555507
$ws.Range("P7").Value = 0.09991617299555505
$ws.Range("Q7").Value = 46.65583382274833
$ws.Range("R7").Value = 419.902504404735
$ws.Range("S7").Value = 0.01029087762044614
$ws.Range("T7").Value = 0.01029087762044613

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.539092333333333
$ws.Range("H8").Value = 7.617277
$ws.Range("I8").Value = 0.1035939158460669
$ws.Range("J8").Value = 0.1035939158460669
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 102.8289443333334
$ws.Range("N8").Value = 308.486833
$ws.Range("O8").Value = 0.5559120396302444
$ws.Range("P8").Value = 0.5559120396302443
$ws.Range("Q8").Value = 261.0921842015268
$ws.Range("R8").Value = 2349.829657813741
$ws.Range("S8").Value = 0.05758910505127092
$ws.Range("T8").Value = 0.05758910505127091

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.539092333333333
$ws.Range("H9").Value = 7.617277
$ws.Range("I9").Value = 0.1035939158460669
$ws.Range("J9").Value = 0.1035939158460669
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 63.66262833333334
$ws.Range("N9").Value = 190.987885
$ws.Range("O9").Value = 0.3441717873742006
$ws.Range("P9").Value = 0.3441717873742006
$ws.Range("Q9").Value = 161.6452915210161
$ws.Range("R9").Value = 1454.807623689145
$ws.Range("S9").Value = 0.03565410317783335
$ws.Range("T9").Value = 0.03565410317783335

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.539092333333333
$ws.Range("H10").Value = 7.617277
$ws.Range("I10").Value = 0.1035939158460669
$ws.Range("J10").Value = 0.1035939158460669
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 18.481835
$ws.Range("N10").Value = 55.445505
$ws.Range("O10").Value = 0.09991617299555507
$ws.Range("P10").Value = 0.09991617299555505
$ws.Range("Q10").Value = 46.92708555443166
$ws.Range("R10").Value = 422.3437699898849
$ws.Range("S10").Value = 0.01035070761696259
$ws.Range("T10").Value = 0.01035070761696259

Write-Output "Applied NATMI recalculated values to Sema3f-Nrp1 sheet"
